$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 4
$ws.Range("F2").Value = 406
$ws.Range("L2").Value = "stimuli/img_juob3.png"
$ws.Range("M2").Value = 79.92105263157895
$ws.Range("N2").Value = 59.78947368421053
$ws.Range("O2").Value = 69.85526315789474
$ws.Range("P2").Value = 38
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = 7
$ws.Range("S2").Value = 7
$ws.Range("T2").Value = 7
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 7

# Row 3
$ws.Range("C3").Value = 4
$ws.Range("F3").Value = 407
$ws.Range("H3").Value = "kitchens"
$ws.Range("I3").Value = "distractor"
$ws.Range("K3").Value = "f"
$ws.Range("L3").Value = "stimuli/img_79b5l.png"
$ws.Range("M3").Value = 72.74285714285715
$ws.Range("N3").Value = 53.31428571428572
$ws.Range("O3").Value = 63.02857142857143
$ws.Range("P3").Value = 35
$ws.Range("Q3").Value = 6
$ws.Range("R3").Value = 6
$ws.Range("S3").Value = 6
$ws.Range("T3").Value = 6
$ws.Range("U3").Value = 6
$ws.Range("V3").Value = 6

# Row 4
$ws.Range("C4").Value = 4
$ws.Range("F4").Value = 408
$ws.Range("H4").Value = "bedrooms"
$ws.Range("I4").Value = "target"
$ws.Range("K4").Value = "j"
$ws.Range("L4").Value = "stimuli/img_rvssl.png"
$ws.Range("M4").Value = 74.25
$ws.Range("N4").Value = 54.33333333333334
$ws.Range("O4").Value = 64.29166666666667
$ws.Range("P4").Value = 36
$ws.Range("V4").Value = 6

# Row 5
$ws.Range("C5").Value = 4
$ws.Range("F5").Value = 409
$ws.Range("L5").Value = "stimuli/img_vh7v8.png"
$ws.Range("M5").Value = 78.70454545454545
$ws.Range("N5").Value = 59.63636363636363
$ws.Range("O5").Value = 69.17045454545455
$ws.Range("Q5").Value = 7
$ws.Range("R5").Value = 7
$ws.Range("S5").Value = 7
$ws.Range("T5").Value = 7
$ws.Range("U5").Value = 7
$ws.Range("V5").Value = 7

# Row 6
$ws.Range("C6").Value = 4
$ws.Range("F6").Value = 410
$ws.Range("H6").Value = "bedrooms"
$ws.Range("I6").Value = "target"
$ws.Range("K6").Value = "j"
$ws.Range("L6").Value = "stimuli/img_wyctg.png"
$ws.Range("M6").Value = 33.44736842105263
$ws.Range("N6").Value = 11.39473684210526
$ws.Range("O6").Value = 22.42105263157895
$ws.Range("P6").Value = 38
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 1
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 1
$ws.Range("U6").Value = 1
$ws.Range("V6").Value = 1

# Row 7
$ws.Range("C7").Value = 4
$ws.Range("F7").Value = 411
$ws.Range("H7").Value = "bedrooms"
$ws.Range("I7").Value = "target"
$ws.Range("K7").Value = "j"
$ws.Range("L7").Value = "stimuli/img_2js6m.png"
$ws.Range("M7").Value = 40.02777777777778
$ws.Range("N7").Value = 20.88888888888889
$ws.Range("O7").Value = 30.45833333333334
$ws.Range("P7").Value = 36
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 2
$ws.Range("S7").Value = 2
$ws.Range("T7").Value = 2
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 2

# Row 8
$ws.Range("C8").Value = 4
$ws.Range("F8").Value = 412
$ws.Range("L8").Value = "stimuli/img_5il0t.png"
$ws.Range("M8").Value = 48.09523809523809
$ws.Range("N8").Value = 30.90476190476191
$ws.Range("O8").Value = 39.5
$ws.Range("P8").Value = 42
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 2
$ws.Range("S8").Value = 2
$ws.Range("T8").Value = 2
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 2

# Row 9
$ws.Range("C9").Value = 4
$ws.Range("F9").Value = 413
$ws.Range("H9").Value = "bedrooms"
$ws.Range("I9").Value = "target"
$ws.Range("K9").Value = "j"
$ws.Range("L9").Value = "stimuli/img_cmyvx.png"
$ws.Range("M9").Value = 64.25
$ws.Range("N9").Value = 40.09375
$ws.Range("O9").Value = 52.171875
$ws.Range("P9").Value = 32
$ws.Range("Q9").Value = 4
$ws.Range("R9").Value = 4
$ws.Range("S9").Value = 4
$ws.Range("T9").Value = 4
$ws.Range("U9").Value = 4
$ws.Range("V9").Value = 4

# Row 10
$ws.Range("C10").Value = 4
$ws.Range("F10").Value = 414
$ws.Range("H10").Value = "kitchens"
$ws.Range("I10").Value = "distractor"
$ws.Range("K10").Value = "f"
$ws.Range("L10").Value = "stimuli/img_e0hwx.png"
$ws.Range("M10").Value = 78.12121212121212
$ws.Range("N10").Value = 55.36363636363637
$ws.Range("O10").Value = 66.74242424242425
$ws.Range("P10").Value = 33
$ws.Range("Q10").Value = 7
$ws.Range("R10").Value = 7
$ws.Range("S10").Value = 7
$ws.Range("T10").Value = 7
$ws.Range("U10").Value = 7
$ws.Range("V10").Value = 7

# Row 11
$ws.Range("C11").Value = 4
$ws.Range("F11").Value = 415
$ws.Range("H11").Value = "bedrooms"
$ws.Range("I11").Value = "target"
$ws.Range("K11").Value = "j"
$ws.Range("L11").Value = "stimuli/img_qihxi.png"
$ws.Range("M11").Value = 76.72222222222223
$ws.Range("N11").Value = 56.33333333333334
$ws.Range("O11").Value = 66.52777777777779
$ws.Range("P11").Value = 36
$ws.Range("Q11").Value = 7
$ws.Range("R11").Value = 7
$ws.Range("S11").Value = 7
$ws.Range("T11").Value = 7
$ws.Range("U11").Value = 7
$ws.Range("V11").Value = 7

# Row 12
$ws.Range("C12").Value = 4
$ws.Range("F12").Value = 416
$ws.Range("H12").Value = "bedrooms"
$ws.Range("I12").Value = "target"
$ws.Range("K12").Value = "j"
$ws.Range("L12").Value = "stimuli/img_le8uf.png"
$ws.Range("M12").Value = 12.88888888888889
$ws.Range("N12").Value = 9.222222222222221
$ws.Range("O12").Value = 11.05555555555556
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = 1
$ws.Range("S12").Value = 1
$ws.Range("T12").Value = 1
$ws.Range("U12").Value = 1
$ws.Range("V12").Value = 1

# Row 13
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 417
$ws.Range("H13").Value = "living_rooms"
$ws.Range("L13").Value = "stimuli/img_eh0no.png"
$ws.Range("M13").Value = 53.66666666666666
$ws.Range("N13").Value = 36.02564102564103
$ws.Range("O13").Value = 44.84615384615385
$ws.Range("P13").Value = 39
$ws.Range("Q13").Value = 3
$ws.Range("R13").Value = 3
$ws.Range("S13").Value = 3
$ws.Range("U13").Value = 3

# Row 14
$ws.Range("C14").Value = 4
$ws.Range("F14").Value = 418
$ws.Range("H14").Value = "bedrooms"
$ws.Range("I14").Value = "target"
$ws.Range("K14").Value = "j"
$ws.Range("L14").Value = "stimuli/img_24rt2.png"
$ws.Range("M14").Value = 55.26829268292683
$ws.Range("N14").Value = 34.19512195121951
$ws.Range("O14").Value = 44.73170731707317
$ws.Range("P14").Value = 41
$ws.Range("Q14").Value = 3
$ws.Range("R14").Value = 3
$ws.Range("S14").Value = 3
$ws.Range("T14").Value = 4
$ws.Range("U14").Value = 4
$ws.Range("V14").Value = 3

# Row 15
$ws.Range("C15").Value = 4
$ws.Range("F15").Value = 419
$ws.Range("H15").Value = "bedrooms"
$ws.Range("I15").Value = "target"
$ws.Range("K15").Value = "j"
$ws.Range("L15").Value = "stimuli/img_cogrz.png"
$ws.Range("M15").Value = 60.5
$ws.Range("N15").Value = 39.71428571428572
$ws.Range("O15").Value = 50.10714285714286
$ws.Range("P15").Value = 42
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = 3
$ws.Range("S15").Value = 3
$ws.Range("T15").Value = 3
$ws.Range("U15").Value = 3
$ws.Range("V15").Value = 3

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("F16").Value = 420
$ws.Range("H16").Value = "kitchens"
$ws.Range("L16").Value = "stimuli/img_0mhms.png"
$ws.Range("M16").Value = 78
$ws.Range("N16").Value = 55.68571428571428
$ws.Range("O16").Value = 66.84285714285714
$ws.Range("P16").Value = 35
$ws.Range("Q16").Value = 7
$ws.Range("R16").Value = 7
$ws.Range("S16").Value = 7
$ws.Range("T16").Value = 7
$ws.Range("U16").Value = 7
$ws.Range("V16").Value = 7

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("F17").Value = 421
$ws.Range("H17").Value = "kitchens"
$ws.Range("L17").Value = "stimuli/img_eppte.png"
$ws.Range("M17").Value = 78.42424242424242
$ws.Range("N17").Value = 57.03030303030303
$ws.Range("O17").Value = 67.72727272727272
$ws.Range("P17").Value = 33
$ws.Range("Q17").Value = 7
$ws.Range("R17").Value = 7
$ws.Range("S17").Value = 7
$ws.Range("T17").Value = 7
$ws.Range("U17").Value = 7
$ws.Range("V17").Value = 7

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("F18").Value = 422
$ws.Range("H18").Value = "kitchens"
$ws.Range("L18").Value = "stimuli/img_xguy9.png"
$ws.Range("M18").Value = 78.21621621621621
$ws.Range("N18").Value = 57.24324324324324
$ws.Range("O18").Value = 67.72972972972973
$ws.Range("P18").Value = 37
$ws.Range("Q18").Value = 7
$ws.Range("R18").Value = 7
$ws.Range("S18").Value = 7
$ws.Range("T18").Value = 7
$ws.Range("U18").Value = 7
$ws.Range("V18").Value = 7

# Row 19
$ws.Range("C19").Value = 4
$ws.Range("F19").Value = 423
$ws.Range("H19").Value = "kitchens"
$ws.Range("I19").Value = "distractor"
$ws.Range("K19").Value = "f"
$ws.Range("L19").Value = "stimuli/img_lpj57.png"
$ws.Range("M19").Value = 74.77777777777777
$ws.Range("N19").Value = 54.44444444444444
$ws.Range("O19").Value = 64.61111111111111
$ws.Range("P19").Value = 27

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("F20").Value = 424
$ws.Range("H20").Value = "living_rooms"
$ws.Range("I20").Value = "distractor"
$ws.Range("K20").Value = "f"
$ws.Range("L20").Value = "stimuli/img_xpco9.png"
$ws.Range("M20").Value = 81.55555555555556
$ws.Range("N20").Value = 64.68888888888888
$ws.Range("O20").Value = 73.12222222222222
$ws.Range("P20").Value = 45
$ws.Range("Q20").Value = 8
$ws.Range("R20").Value = 8
$ws.Range("S20").Value = 8
$ws.Range("T20").Value = 8
$ws.Range("U20").Value = 8
$ws.Range("V20").Value = 8

# Row 21
$ws.Range("C21").Value = 4
$ws.Range("F21").Value = 425
$ws.Range("L21").Value = "stimuli/img_x0u5z.png"
$ws.Range("M21").Value = 92
$ws.Range("N21").Value = 78.16216216216216
$ws.Range("O21").Value = 85.08108108108108
$ws.Range("P21").Value = 37
$ws.Range("Q21").Value = 10
$ws.Range("R21").Value = 10
$ws.Range("S21").Value = 10
$ws.Range("T21").Value = 10
$ws.Range("U21").Value = 10
$ws.Range("V21").Value = 10

# Row 22
$ws.Range("C22").Value = 4
$ws.Range("F22").Value = 426
$ws.Range("L22").Value = "stimuli/img_3h4c9.png"
$ws.Range("M22").Value = 85.47619047619048
$ws.Range("N22").Value = 67.26190476190476
$ws.Range("O22").Value = 76.36904761904762
$ws.Range("P22").Value = 42
$ws.Range("Q22").Value = 9
$ws.Range("R22").Value = 9
$ws.Range("S22").Value = 9
$ws.Range("T22").Value = 9
$ws.Range("U22").Value = 9
$ws.Range("V22").Value = 9

# Row 23
$ws.Range("C23").Value = 4
$ws.Range("F23").Value = 427
$ws.Range("H23").Value = "bedrooms"
$ws.Range("I23").Value = "target"
$ws.Range("K23").Value = "j"
$ws.Range("L23").Value = "stimuli/img_72fmj.png"
$ws.Range("M23").Value = 53.87179487179487
$ws.Range("N23").Value = 36.02564102564103
$ws.Range("O23").Value = 44.94871794871795
$ws.Range("P23").Value = 39
$ws.Range("Q23").Value = 3
$ws.Range("R23").Value = 3
$ws.Range("S23").Value = 3
$ws.Range("T23").Value = 3
$ws.Range("U23").Value = 3
$ws.Range("V23").Value = 3

# Row 24
$ws.Range("C24").Value = 4
$ws.Range("F24").Value = 428
$ws.Range("H24").Value = "living_rooms"
$ws.Range("I24").Value = "distractor"
$ws.Range("K24").Value = "f"
$ws.Range("L24").Value = "stimuli/img_3m61b.png"
$ws.Range("M24").Value = 81.97619047619048
$ws.Range("N24").Value = 63.23809523809524
$ws.Range("O24").Value = 72.60714285714286
$ws.Range("P24").Value = 42
$ws.Range("Q24").Value = 8
$ws.Range("R24").Value = 8
$ws.Range("S24").Value = 8
$ws.Range("T24").Value = 8
$ws.Range("U24").Value = 8
$ws.Range("V24").Value = 8

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("F25").Value = 429
$ws.Range("H25").Value = "living_rooms"
$ws.Range("L25").Value = "stimuli/img_3jnt7.png"
$ws.Range("M25").Value = 49.52272727272727
$ws.Range("N25").Value = 35.25
$ws.Range("O25").Value = 42.38636363636364
$ws.Range("P25").Value = 44
$ws.Range("Q25").Value = 3
$ws.Range("R25").Value = 3
$ws.Range("S25").Value = 3
$ws.Range("T25").Value = 3
$ws.Range("U25").Value = 3

# Row 26
$ws.Range("C26").Value = 4
$ws.Range("F26").Value = 430
$ws.Range("H26").Value = "kitchens"
$ws.Range("I26").Value = "distractor"
$ws.Range("K26").Value = "f"
$ws.Range("L26").Value = "stimuli/img_60242.png"
$ws.Range("M26").Value = 78.33333333333333
$ws.Range("N26").Value = 57.57575757575758
$ws.Range("O26").Value = 67.95454545454545
$ws.Range("P26").Value = 33
$ws.Range("Q26").Value = 7
$ws.Range("R26").Value = 7
$ws.Range("S26").Value = 7
$ws.Range("T26").Value = 7
$ws.Range("U26").Value = 7
$ws.Range("V26").Value = 7

# Row 27
$ws.Range("C27").Value = 4
$ws.Range("F27").Value = 431
$ws.Range("H27").Value = "living_rooms"
$ws.Range("I27").Value = "distractor"
$ws.Range("K27").Value = "f"
$ws.Range("L27").Value = "stimuli/img_89dvt.png"
$ws.Range("M27").Value = 81.09756097560975
$ws.Range("N27").Value = 64.6829268292683
$ws.Range("O27").Value = 72.89024390243902
$ws.Range("P27").Value = 41
$ws.Range("Q27").Value = 8
$ws.Range("R27").Value = 8
$ws.Range("S27").Value = 8
$ws.Range("T27").Value = 8
$ws.Range("U27").Value = 8
$ws.Range("V27").Value = 8

# Row 28
$ws.Range("C28").Value = 4
$ws.Range("F28").Value = 432
$ws.Range("H28").Value = "kitchens"
$ws.Range("I28").Value = "distractor"
$ws.Range("K28").Value = "f"
$ws.Range("L28").Value = "stimuli/img_cxpff.png"
$ws.Range("M28").Value = 74.92307692307692
$ws.Range("N28").Value = 53.28205128205128
$ws.Range("O28").Value = 64.1025641025641
$ws.Range("P28").Value = 39
$ws.Range("Q28").Value = 6
$ws.Range("R28").Value = 6
$ws.Range("S28").Value = 6
$ws.Range("T28").Value = 6
$ws.Range("U28").Value = 6
$ws.Range("V28").Value = 6

# Row 29
$ws.Range("C29").Value = 4
$ws.Range("F29").Value = 433
$ws.Range("L29").Value = "stimuli/img_oou46.png"
$ws.Range("M29").Value = 75.70270270270271
$ws.Range("N29").Value = 54.86486486486486
$ws.Range("O29").Value = 65.28378378378379
$ws.Range("P29").Value = 37
$ws.Range("Q29").Value = 6
$ws.Range("R29").Value = 6
$ws.Range("S29").Value = 6
$ws.Range("T29").Value = 6
$ws.Range("U29").Value = 6
$ws.Range("V29").Value = 6

# Row 30
$ws.Range("C30").Value = 4
$ws.Range("F30").Value = 434
$ws.Range("L30").Value = "stimuli/img_uxxo0.png"
$ws.Range("M30").Value = 71.74418604651163
$ws.Range("N30").Value = 48.44186046511628
$ws.Range("O30").Value = 60.09302325581395
$ws.Range("P30").Value = 43
$ws.Range("Q30").Value = 5
$ws.Range("R30").Value = 5
$ws.Range("S30").Value = 5
$ws.Range("T30").Value = 5
$ws.Range("U30").Value = 5
$ws.Range("V30").Value = 5

# Row 31
$ws.Range("C31").Value = 4
$ws.Range("F31").Value = 435
$ws.Range("H31").Value = "bedrooms"
$ws.Range("I31").Value = "target"
$ws.Range("K31").Value = "j"
$ws.Range("L31").Value = "stimuli/img_fnu4h.png"
$ws.Range("M31").Value = 85.87179487179488
$ws.Range("N31").Value = 70.71794871794872
$ws.Range("O31").Value = 78.2948717948718
$ws.Range("P31").Value = 39
$ws.Range("Q31").Value = 9
$ws.Range("R31").Value = 9
$ws.Range("S31").Value = 9
$ws.Range("T31").Value = 9
$ws.Range("U31").Value = 9
$ws.Range("V31").Value = 9

# Row 32
$ws.Range("C32").Value = 4
$ws.Range("F32").Value = 436
$ws.Range("L32").Value = "stimuli/img_1vq1v.png"
$ws.Range("M32").Value = 69.42857142857143
$ws.Range("N32").Value = 46.59523809523809
$ws.Range("O32").Value = 58.01190476190476
$ws.Range("P32").Value = 42
$ws.Range("Q32").Value = 5
$ws.Range("R32").Value = 5
$ws.Range("S32").Value = 5
$ws.Range("T32").Value = 5
$ws.Range("U32").Value = 5
$ws.Range("V32").Value = 5

# Row 33
$ws.Range("C33").Value = 4
$ws.Range("F33").Value = 437
$ws.Range("H33").Value = "living_rooms"
$ws.Range("I33").Value = "distractor"
$ws.Range("K33").Value = "f"
$ws.Range("L33").Value = "stimuli/img_cehin.png"
$ws.Range("M33").Value = 78.86363636363636
$ws.Range("N33").Value = 60.02272727272727
$ws.Range("O33").Value = 69.44318181818181
$ws.Range("P33").Value = 44
$ws.Range("Q33").Value = 7
$ws.Range("R33").Value = 7
$ws.Range("S33").Value = 7
$ws.Range("T33").Value = 7
$ws.Range("U33").Value = 7
$ws.Range("V33").Value = 7

# Row 34
$ws.Range("C34").Value = 4
$ws.Range("F34").Value = 438
$ws.Range("L34").Value = "stimuli/img_f4jxo.png"
$ws.Range("M34").Value = 82.91666666666667
$ws.Range("N34").Value = 65.52777777777777
$ws.Range("O34").Value = 74.22222222222223
$ws.Range("P34").Value = 36
$ws.Range("Q34").Value = 8
$ws.Range("R34").Value = 8
$ws.Range("S34").Value = 8
$ws.Range("T34").Value = 8
$ws.Range("U34").Value = 8
$ws.Range("V34").Value = 8

# Row 35
$ws.Range("C35").Value = 4
$ws.Range("F35").Value = 439
$ws.Range("L35").Value = "stimuli/img_a9acb.png"
$ws.Range("M35").Value = 77.11428571428571
$ws.Range("N35").Value = 58.42857142857143
$ws.Range("O35").Value = 67.77142857142857
$ws.Range("P35").Value = 35
$ws.Range("Q35").Value = 7
$ws.Range("R35").Value = 7
$ws.Range("S35").Value = 7
$ws.Range("T35").Value = 7
$ws.Range("U35").Value = 7
$ws.Range("V35").Value = 7

# Row 36
$ws.Range("C36").Value = 4
$ws.Range("F36").Value = 440
$ws.Range("H36").Value = "bedrooms"
$ws.Range("I36").Value = "target"
$ws.Range("K36").Value = "j"
$ws.Range("L36").Value = "stimuli/img_t2ioc.png"
$ws.Range("M36").Value = 88.1891891891892
$ws.Range("N36").Value = 74.05405405405405
$ws.Range("O36").Value = 81.12162162162161
$ws.Range("P36").Value = 37
$ws.Range("Q36").Value = 10
$ws.Range("R36").Value = 10
$ws.Range("S36").Value = 10
$ws.Range("T36").Value = 10
$ws.Range("U36").Value = 10
$ws.Range("V36").Value = 10

# Row 37
$ws.Range("C37").Value = 4
$ws.Range("F37").Value = 441
$ws.Range("H37").Value = "bedrooms"
$ws.Range("I37").Value = "target"
$ws.Range("K37").Value = "j"
$ws.Range("L37").Value = "stimuli/img_jp28n.png"
$ws.Range("M37").Value = 65.02564102564102
$ws.Range("N37").Value = 44.97435897435897
$ws.Range("O37").Value = 55
$ws.Range("P37").Value = 39
$ws.Range("V37").Value = 5

# Row 38
$ws.Range("C38").Value = 4
$ws.Range("F38").Value = 442
$ws.Range("H38").Value = "bedrooms"
$ws.Range("I38").Value = "target"
$ws.Range("K38").Value = "j"
$ws.Range("L38").Value = "stimuli/img_kugyw.png"
$ws.Range("M38").Value = 74.25
$ws.Range("N38").Value = 54.10714285714285
$ws.Range("O38").Value = 64.17857142857143
$ws.Range("P38").Value = 28
$ws.Range("Q38").Value = 6
$ws.Range("R38").Value = 6
$ws.Range("S38").Value = 6
$ws.Range("T38").Value = 6
$ws.Range("U38").Value = 6
$ws.Range("V38").Value = 6

# Row 39
$ws.Range("C39").Value = 4
$ws.Range("F39").Value = 443
$ws.Range("L39").Value = "stimuli/img_e26ut.png"
$ws.Range("M39").Value = 81.07692307692308
$ws.Range("N39").Value = 61.28205128205128
$ws.Range("O39").Value = 71.17948717948718
$ws.Range("P39").Value = 39
$ws.Range("R39").Value = 8
$ws.Range("S39").Value = 8
$ws.Range("T39").Value = 8
$ws.Range("U39").Value = 8
$ws.Range("V39").Value = 8

# Row 40
$ws.Range("C40").Value = 4
$ws.Range("F40").Value = 444
$ws.Range("H40").Value = "bedrooms"
$ws.Range("I40").Value = "target"
$ws.Range("K40").Value = "j"
$ws.Range("L40").Value = "stimuli/img_7lz7m.png"
$ws.Range("M40").Value = 51.5531914893617
$ws.Range("N40").Value = 32.87234042553192
$ws.Range("O40").Value = 42.21276595744681
$ws.Range("P40").Value = 47
$ws.Range("Q40").Value = 3
$ws.Range("R40").Value = 3
$ws.Range("S40").Value = 3
$ws.Range("T40").Value = 3
$ws.Range("U40").Value = 3
$ws.Range("V40").Value = 3

# Row 41
$ws.Range("C41").Value = 4
$ws.Range("F41").Value = 445
$ws.Range("H41").Value = "living_rooms"
$ws.Range("I41").Value = "distractor"
$ws.Range("K41").Value = "f"
$ws.Range("L41").Value = "stimuli/img_emh91.png"
$ws.Range("M41").Value = 82.06666666666666
$ws.Range("N41").Value = 63.33333333333334
$ws.Range("O41").Value = 72.7
$ws.Range("P41").Value = 45
$ws.Range("Q41").Value = 8
$ws.Range("R41").Value = 8
$ws.Range("S41").Value = 8
$ws.Range("T41").Value = 8
$ws.Range("U41").Value = 8
$ws.Range("V41").Value = 8
